$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C3: was a text label ("2 420,5*"), now a numeric value styled like the other cells.
$c3 = $ws.Range("C3")
$c3.Value = 2420.5
$c3.NumberFormat = "#,##0.00"

# Row 1 height change (106 -> 100)
$ws.Rows.Item(1).RowHeight = 100

# Selection / view change
$ws.Range("D7").Select() | Out-Null
$excel.ActiveWindow.Zoom = 159

Write-Host "Done"
